$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44568
$ws.Range("K2").Value = 'Santina'
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = '$/bandeja 12 kilos'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1292
$ws.Range("T2").Value = 12

# Row 4
$ws.Range("D4").Value2 = 44537
$ws.Range("K4").Value = 'Brooks'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 29000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 29500
$ws.Range("Q4").Value = '$/caja 20 kilos'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1475
$ws.Range("T4").Value = 20

# Row 5
$ws.Range("D5").Value2 = 44175
$ws.Range("K5").Value = 'Rainier'
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 26000
$ws.Range("P5").Value = 25500
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 1417
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value2 = 44161
$ws.Range("K6").Value = 'Bing'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 39000
$ws.Range("O6").Value = 40000
$ws.Range("P6").Value = 39500
$ws.Range("Q6").Value = '$/caja 20 kilos'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 1975
$ws.Range("T6").Value = 20

# Row 8
$ws.Range("D8").Value2 = 44208
$ws.Range("K8").Value = 'Lapins'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 10500
$ws.Range("O8").Value = 11000
$ws.Range("P8").Value = 10750
$ws.Range("Q8").Value = '$/bandeja 12 kilos'
$ws.Range("R8").Value = 'Provincia de Curicó'
$ws.Range("S8").Value = 896
$ws.Range("T8").Value = 12

# Row 9
$ws.Range("D9").Value2 = 44229
$ws.Range("K9").Value = 'Santina'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 6500
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 6750
$ws.Range("Q9").Value = '$/bandeja 5 kilos'
$ws.Range("R9").Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 1350
$ws.Range("T9").Value = 5

# Row 10
$ws.Range("D10").Value2 = 44532
$ws.Range("K10").Value = 'Brooks'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 400
$ws.Range("N10").Value = 27000
$ws.Range("O10").Value = 28000
$ws.Range("P10").Value = 27500
$ws.Range("Q10").Value = '$/bandeja 12 kilos'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 2292
$ws.Range("T10").Value = 12

# Row 11
$ws.Range("D11").Value2 = 44571
$ws.Range("K11").Value = 'Brooks'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 400
$ws.Range("N11").Value = 8500
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 8750
$ws.Range("Q11").Value = '$/bandeja 10 kilos'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 875
$ws.Range("T11").Value = 10

